# Weekly fruit/vegetable price update:
# Insert a new daily observation at row 22 (pushing the existing rows 22-37
# down to 23-38), then populate the new row with the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 22:37 down to 23:38 to make room for the new record.
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with the new weekly data point.
$ws.Range("A22").Value = 4
$ws.Range("B22").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C22").Value = "Los Lagos"
$ws.Range("D22").Value = 44460
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = 100112026
$ws.Range("G22").Value = "Haba"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 15000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 15000
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Provincia de Limarí"
$ws.Range("P22").Value = 600
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
